$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header strings for columns L through S (row 1)
$ws.Range("L1").Value = "hzj-混合调节_20170516_152754_ASIC_EEG"
$ws.Range("M1").Value = "hzj-混合调节_20170518_134207_ASIC_EEG"
$ws.Range("N1").Value = "hzj-混合调节_20170519_135415_ASIC_EEG"
$ws.Range("O1").Value = "zyx-混合调节_20170516_111228_ASIC_EEG"
$ws.Range("P1").Value = "zyx-混合调节_20170517_110944_ASIC_EEG"
$ws.Range("Q1").Value = "zyx-混合调节_20170518_112337_ASIC_EEG"
$ws.Range("R1").Value = "zyx-混合调节_20170519_124954_ASIC_EEG"
$ws.Range("S1").Value = "zyx-混合调节_20170522_111557_ASIC_EEG"

# Row 2 numeric values
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.95876288659793818
$ws.Range("N2").Value = 0.94756554307116103
$ws.Range("O2").Value = 0.96153846153846156
$ws.Range("P2").Value = 0.94871794871794868
$ws.Range("Q2").Value = 0.93527508090614886
$ws.Range("R2").Value = 0.97419354838709682
$ws.Range("S2").Value = 0.94660194174757284

# Row 3 numeric values
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 0.97741935483870968
$ws.Range("N3").Value = 1
$ws.Range("O3").Value = 0.98928571428571432
$ws.Range("P3").Value = 0.97206703910614523
$ws.Range("Q3").Value = 0.97938144329896903
$ws.Range("R3").Value = 0.96308724832214765
$ws.Range("S3").Value = 0.95238095238095233

# Update selection to reflect the new used range
$ws.Range("A1:S3").Select()
